$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row: E1 "Stuff2" header becomes "Scholar" ---
$ws.Range("E1").Value = "Scholar"

# --- Spell learning (column B) progression changed from shared FLOOR.MATH formula to static values ---
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 3
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 3
$ws.Range("B12").Value = 4
$ws.Range("B13").Value = 4
$ws.Range("B14").Value = 4
$ws.Range("B15").Value = 5
$ws.Range("B16").Value = 5
$ws.Range("B17").Value = 5
$ws.Range("B18").Value = 6
$ws.Range("B19").Value = 6
$ws.Range("B20").Value = 6
$ws.Range("B21").Value = 7

# --- Clear the old row-2 "Stuff0/Stuff1/Stuff2" placeholder cells ---
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# --- Scholar archetype content, written in the order the new entries were authored ---
$ws.Range("D4").Value = "Instructive Aura"
$ws.Range("D7").Value = "Guiding Hand"
$ws.Range("D11").Value = "Guiding Hand II"
$ws.Range("C2").Value = "Intellectual Insight"
$ws.Range("C3").Value = "Arcane Affinity"
$ws.Range("D6").Value = "Repository of Knowledge"
$ws.Range("E8").Value = "Rapid Assimilation"
$ws.Range("E5").Value = "Innovative Mind"
$ws.Range("E4").Value = "Research Training"

# --- Column widths: split the old merged C:D width into distinct C, D, E widths ---
$ws.Columns.Item(4).ColumnWidth = 32.09
$ws.Columns.Item(5).ColumnWidth = 22.25

# --- Selection moved to E10 ---
$ws.Range("E10").Select()
